# Auto-generated script applying scheduled market-data refresh to the
# currentAveragePrice / LevePrice / LeveProfit columns across all sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 236.73077
$ws.Range("I9").Value = 393.2857
$ws.Range("J9").Value = 179.05263
$ws.Range("K9").Value = 393.2857
$ws.Range("L9").Value = 179.05263
$ws.Range("M9").Value = -224.2857
$ws.Range("N9").Value = -517.05263
$ws.Range("H11").Value = 68.13333
$ws.Range("I11").Value = 68.13333
$ws.Range("K11").Value = 68.13333
$ws.Range("M11").Value = 71.86667
$ws.Range("H106").Value = 8284.5
$ws.Range("I106").Value = 1791.3
$ws.Range("J106").Value = 19106.5
$ws.Range("K106").Value = 1791.3
$ws.Range("L106").Value = 19106.5
$ws.Range("M106").Value = -1160.3
$ws.Range("N106").Value = -20368.5
$ws.Range("H107").Value = 887.5
$ws.Range("I107").Value = 887.5
$ws.Range("K107").Value = 887.5
$ws.Range("M107").Value = 1032.5
$ws.Range("H132").Value = 3414.6365
$ws.Range("I132").Value = 2889.8333
$ws.Range("J132").Value = 5776.25
$ws.Range("K132").Value = 8669.499899999999
$ws.Range("L132").Value = 17328.75
$ws.Range("M132").Value = -6139.499899999999
$ws.Range("N132").Value = -22388.75
$ws.Range("H138").Value = 3076.1143
$ws.Range("I138").Value = 3952.5
$ws.Range("J138").Value = 3023
$ws.Range("K138").Value = 11857.5
$ws.Range("L138").Value = 9069
$ws.Range("M138").Value = -6717.5
$ws.Range("N138").Value = -19349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3286.2222
$ws.Range("I2").Value = 2947.625
$ws.Range("J2").Value = 5995
$ws.Range("K2").Value = 2947.625
$ws.Range("L2").Value = 5995
$ws.Range("M2").Value = -2834.625
$ws.Range("N2").Value = -6221
$ws.Range("H43").Value = 19000
$ws.Range("J43").Value = 19000
$ws.Range("L43").Value = 19000
$ws.Range("N43").Value = -19626
$ws.Range("H61").Value = 1892.5333
$ws.Range("I61").Value = 1688.6923
$ws.Range("J61").Value = 3217.5
$ws.Range("K61").Value = 1688.6923
$ws.Range("L61").Value = 3217.5
$ws.Range("M61").Value = -1476.6923
$ws.Range("N61").Value = -3641.5
$ws.Range("H74").Value = 1311.2354
$ws.Range("I74").Value = 1399.3572
$ws.Range("K74").Value = 1399.3572
$ws.Range("M74").Value = -525.3571999999999
$ws.Range("H77").Value = 1311.2354
$ws.Range("I77").Value = 1399.3572
$ws.Range("K77").Value = 6996.786
$ws.Range("M77").Value = -2628.786
$ws.Range("H116").Value = 3286.2222
$ws.Range("I116").Value = 2947.625
$ws.Range("J116").Value = 5995
$ws.Range("K116").Value = 2947.625
$ws.Range("L116").Value = 5995
$ws.Range("M116").Value = -653.625
$ws.Range("N116").Value = -10583
$ws.Range("H132").Value = 12886.421
$ws.Range("I132").Value = 14959.125
$ws.Range("J132").Value = 1832
$ws.Range("K132").Value = 44877.375
$ws.Range("L132").Value = 5496
$ws.Range("M132").Value = -42347.375
$ws.Range("N132").Value = -10556
$ws.Range("H136").Value = 1892.5333
$ws.Range("I136").Value = 1688.6923
$ws.Range("J136").Value = 3217.5
$ws.Range("K136").Value = 5066.0769
$ws.Range("L136").Value = 9652.5
$ws.Range("M136").Value = -2516.0769
$ws.Range("N136").Value = -14752.5
$ws.Range("H139").Value = 124901.664
$ws.Range("J139").Value = 93355
$ws.Range("L139").Value = 93355
$ws.Range("N139").Value = -103635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3286.2222
$ws.Range("I3").Value = 2947.625
$ws.Range("J3").Value = 5995
$ws.Range("K3").Value = 2947.625
$ws.Range("L3").Value = 5995
$ws.Range("M3").Value = -2833.625
$ws.Range("N3").Value = -6223
$ws.Range("H64").Value = 924
$ws.Range("I64").Value = 90
$ws.Range("J64").Value = 1202
$ws.Range("K64").Value = 90
$ws.Range("L64").Value = 1202
$ws.Range("M64").Value = 135
$ws.Range("N64").Value = -1652
$ws.Range("H67").Value = 924
$ws.Range("I67").Value = 90
$ws.Range("J67").Value = 1202
$ws.Range("K67").Value = 90
$ws.Range("L67").Value = 1202
$ws.Range("M67").Value = 690
$ws.Range("N67").Value = -2762
$ws.Range("H86").Value = 1784.5333
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 1652.5714
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 1652.5714
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -3898.5714
$ws.Range("H89").Value = 1784.5333
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 1652.5714
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 8262.857
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -19494.857
$ws.Range("H105").Value = 5887.6113
$ws.Range("I105").Value = 6477.8
$ws.Range("J105").Value = 2936.6667
$ws.Range("K105").Value = 6477.8
$ws.Range("L105").Value = 2936.6667
$ws.Range("M105").Value = -4730.8
$ws.Range("N105").Value = -6430.6667
$ws.Range("H134").Value = 1881.7142
$ws.Range("I134").Value = 1650.3793
$ws.Range("K134").Value = 4951.1379
$ws.Range("M134").Value = -2416.1379
$ws.Range("H138").Value = 65000
$ws.Range("I138").Value = 65000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 65000
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("M138").Value = -59860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1059.1072
$ws.Range("I16").Value = 1062.037
$ws.Range("J16").Value = 980
$ws.Range("K16").Value = 1062.037
$ws.Range("L16").Value = 980
$ws.Range("M16").Value = -775.037
$ws.Range("N16").Value = -1554
$ws.Range("H19").Value = 492.4
$ws.Range("I19").Value = 309.25
$ws.Range("J19").Value = 1225
$ws.Range("K19").Value = 309.25
$ws.Range("L19").Value = 1225
$ws.Range("M19").Value = -139.25
$ws.Range("N19").Value = -1565
$ws.Range("H24").Value = 492.4
$ws.Range("I24").Value = 309.25
$ws.Range("J24").Value = 1225
$ws.Range("K24").Value = 309.25
$ws.Range("L24").Value = 1225
$ws.Range("M24").Value = -139.25
$ws.Range("N24").Value = -1565
$ws.Range("H31").Value = 2585.2
$ws.Range("I31").Value = 2264.7856
$ws.Range("J31").Value = 2993
$ws.Range("K31").Value = 2264.7856
$ws.Range("L31").Value = 2993
$ws.Range("M31").Value = -1969.7856
$ws.Range("N31").Value = -3583
$ws.Range("H34").Value = 2585.2
$ws.Range("I34").Value = 2264.7856
$ws.Range("J34").Value = 2993
$ws.Range("K34").Value = 2264.7856
$ws.Range("L34").Value = 2993
$ws.Range("M34").Value = -2062.7856
$ws.Range("N34").Value = -3397
$ws.Range("H58").Value = 2108
$ws.Range("I58").Value = 1853
$ws.Range("J58").Value = 2312
$ws.Range("K58").Value = 1853
$ws.Range("L58").Value = 2312
$ws.Range("M58").Value = -1650
$ws.Range("N58").Value = -2718
$ws.Range("H99").Value = 18421.105
$ws.Range("I99").Value = 32749.875
$ws.Range("K99").Value = 32749.875
$ws.Range("M99").Value = -31251.875
$ws.Range("H105").Value = 969.1613
$ws.Range("I105").Value = 916.875
$ws.Range("J105").Value = 1148.4286
$ws.Range("K105").Value = 916.875
$ws.Range("L105").Value = 1148.4286
$ws.Range("M105").Value = 830.125
$ws.Range("N105").Value = -4642.4286
$ws.Range("H113").Value = 1059.1072
$ws.Range("I113").Value = 1062.037
$ws.Range("J113").Value = 980
$ws.Range("K113").Value = 1062.037
$ws.Range("L113").Value = 980
$ws.Range("M113").Value = 1107.963
$ws.Range("N113").Value = -5320
$ws.Range("H126").Value = 18421.105
$ws.Range("I126").Value = 32749.875
$ws.Range("K126").Value = 98249.625
$ws.Range("M126").Value = -95779.625
$ws.Range("H132").Value = 2826.0344
$ws.Range("I132").Value = 2753.7083
$ws.Range("J132").Value = 3173.2
$ws.Range("K132").Value = 8261.124899999999
$ws.Range("L132").Value = 9519.599999999999
$ws.Range("M132").Value = -5731.124899999999
$ws.Range("N132").Value = -14579.6
$ws.Range("H134").Value = 2370.973
$ws.Range("I134").Value = 2309.8857
$ws.Range("K134").Value = 6929.657099999999
$ws.Range("M134").Value = -4394.657099999999
$ws.Range("H136").Value = 2108
$ws.Range("I136").Value = 1853
$ws.Range("J136").Value = 2312
$ws.Range("K136").Value = 5559
$ws.Range("L136").Value = 6936
$ws.Range("M136").Value = -3009
$ws.Range("N136").Value = -12036

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1035.7142
$ws.Range("I107").Value = 1324.75
$ws.Range("J107").Value = 987.5417
$ws.Range("K107").Value = 3974.25
$ws.Range("L107").Value = 2962.6251
$ws.Range("M107").Value = -2054.25
$ws.Range("N107").Value = -6802.6251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2726.963
$ws.Range("I132").Value = 2908.5217
$ws.Range("J132").Value = 1683
$ws.Range("K132").Value = 8725.5651
$ws.Range("L132").Value = 5049
$ws.Range("M132").Value = -6195.5651
$ws.Range("N132").Value = -10109

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5314.353
$ws.Range("I100").Value = 4653.75
$ws.Range("J100").Value = 6899.8
$ws.Range("K100").Value = 4653.75
$ws.Range("L100").Value = 6899.8
$ws.Range("M100").Value = -4112.75
$ws.Range("N100").Value = -7981.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21745202
$ws.Range("I107").Value = 7536
$ws.Range("J107").Value = 83335256
$ws.Range("K107").Value = 22608
$ws.Range("L107").Value = 250005768
$ws.Range("M107").Value = -20688
$ws.Range("N107").Value = -250009608
$ws.Range("H132").Value = 6832.4443
$ws.Range("I132").Value = 9058.474
$ws.Range("J132").Value = 1545.625
$ws.Range("K132").Value = 27175.422
$ws.Range("L132").Value = 4636.875
$ws.Range("M132").Value = -24645.422
$ws.Range("N132").Value = -9696.875
